function Set-CellText($ws, $addr, $val) {
    if ($val -match '^[+-]?\d+(\.\d+)?$') {
        $ws.Range($addr).Value = "'" + $val
    } else {
        $ws.Range($addr).Value = $val
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
Set-CellText $ws "D2" "19.947.29"
Set-CellText $ws "E2" "  -8.01%  "

# Row 3
Set-CellText $ws "D3" "1.408.82"
Set-CellText $ws "E3" "  -8.22%  "

# Row 4
Set-CellText $ws "E4" "  +0.01%  "

# Row 5
Set-CellText $ws "E5" "  +0.02%  "

# Row 6
Set-CellText $ws "D6" "272.77"
Set-CellText $ws "E6" "  -5.67%  "

# Row 7
Set-CellText $ws "D7" "0.3690"
Set-CellText $ws "E7" "  -6.19%  "

# Row 8
Set-CellText $ws "D8" "0.3070"
Set-CellText $ws "E8" "  -2.63%  "

# Row 9
Set-CellText $ws "D9" "39.16"
Set-CellText $ws "E9" "  -7.34%  "

# Row 10
Set-CellText $ws "D10" "0.9890"
Set-CellText $ws "E10" "  -5.44%  "

# Row 11
Set-CellText $ws "D11" "0.06534"
Set-CellText $ws "E11" "  -8.80%  "

# Row 12
Set-CellText $ws "D12" "1.002"
Set-CellText $ws "E12" "  +0.02%  "

# Row 13
Set-CellText $ws "D13" "5.310"
Set-CellText $ws "E13" "  -5.25%  "

# Row 14
Set-CellText $ws "B14" "Chainlink"
Set-CellText $ws "C14" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-CellText $ws "D14" "6.129"
Set-CellText $ws "E14" "  -7.01%  "

# Row 15
Set-CellText $ws "B15" "Solana"
Set-CellText $ws "C15" "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-CellText $ws "D15" "16.87"
Set-CellText $ws "E15" "  -8.48%  "

# Row 16
Set-CellText $ws "D16" "1.410.54"
Set-CellText $ws "E16" "  -8.16%  "

# Row 17
Set-CellText $ws "D17" "0.00001004"
Set-CellText $ws "E17" "  -8.41%  "

# Row 18
Set-CellText $ws "D18" "0.05758"
Set-CellText $ws "E18" "  -12.47%  "

# Row 19
Set-CellText $ws "D19" "73.05"
Set-CellText $ws "E19" "  -11.86%  "

# Row 20
Set-CellText $ws "E20" "  +0.05%  "

# Row 21
Set-CellText $ws "D21" "5.574"
Set-CellText $ws "E21" "  -8.80%  "

# Row 22
Set-CellText $ws "D22" "14.33"
Set-CellText $ws "E22" "  -6.55%  "

# Row 23
Set-CellText $ws "D23" "10.77"
Set-CellText $ws "E23" "  -0.65%  "

# Row 24
Set-CellText $ws "D24" "2.274"
Set-CellText $ws "E24" "  -4.57%  "

# Row 25
Set-CellText $ws "D25" "19.949.16"
Set-CellText $ws "E25" "  -8.00%  "

# Row 26
Set-CellText $ws "D26" "2.231"
Set-CellText $ws "E26" "  -4.72%  "

# Row 27
Set-CellText $ws "D27" "138.15"
Set-CellText $ws "E27" "  -5.59%  "

# Row 28
Set-CellText $ws "D28" "16.71"
Set-CellText $ws "E28" "  -8.60%  "

# Row 29
Set-CellText $ws "D29" "1.569.11"
Set-CellText $ws "E29" "  -8.17%  "

# Row 30
Set-CellText $ws "D30" "108.60"
Set-CellText $ws "E30" "  -7.06%  "

# Row 31
Set-CellText $ws "D31" "3.833"
Set-CellText $ws "E31" "  -20.75%  "

# Row 32
Set-CellText $ws "D32" "5.275"
Set-CellText $ws "E32" "  -9.86%  "

# Row 33
Set-CellText $ws "D33" "0.8095"
Set-CellText $ws "E33" "  -15.59%  "

# Row 34
Set-CellText $ws "D34" "0.07669"
Set-CellText $ws "E34" "  -6.00%  "

# Row 35
Set-CellText $ws "D35" "8.427"
Set-CellText $ws "E35" "  -2.52%  "

# Row 36
Set-CellText $ws "D36" "0.05755"
Set-CellText $ws "E36" "  -4.90%  "

# Row 37
Set-CellText $ws "B37" "Frax"
Set-CellText $ws "C37" "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-CellText $ws "D37" "1.001"
Set-CellText $ws "E37" "  +0.02%  "

# Row 38
Set-CellText $ws "B38" "InternetComputer(DFINITY)"
Set-CellText $ws "C38" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-CellText $ws "D38" "4.767"
Set-CellText $ws "E38" "  -6.32%  "

# Row 39
Set-CellText $ws "D39" "0.1927"
Set-CellText $ws "E39" "  -4.56%  "

# Row 40
Set-CellText $ws "D40" "0.02032"
Set-CellText $ws "E40" "  -7.36%  "

# Row 41
Set-CellText $ws "D41" "10.22"
Set-CellText $ws "E41" "  -4.10%  "

# Row 42
Set-CellText $ws "D42" "1.054"
Set-CellText $ws "E42" "  -10.68%  "

# Row 43
Set-CellText $ws "D43" "1.275"
Set-CellText $ws "E43" "  -11.29%  "

# Row 44
Set-CellText $ws "D44" "0.5264"
Set-CellText $ws "E44" "  -7.75%  "

# Row 45
Set-CellText $ws "D45" "3.519"
Set-CellText $ws "E45" "  -5.49%  "

# Row 46
Set-CellText $ws "D46" "12.18"
Set-CellText $ws "E46" "  -6.58%  "

# Row 47
Set-CellText $ws "D47" "0.5084"
Set-CellText $ws "E47" "  -6.92%  "

# Row 48
Set-CellText $ws "D48" "1.795"
Set-CellText $ws "E48" "  -3.24%  "

# Row 49
Set-CellText $ws "D49" "110.49"
Set-CellText $ws "E49" "  -4.73%  "

# Row 50
Set-CellText $ws "D50" "1.032"
Set-CellText $ws "E50" "  -11.27%  "

# Row 51
Set-CellText $ws "E51" "  +0.03%  "
